# Apply updates to column F ("dSF") for specific rows as described in the
# commit: "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -4
    5  = 3
    9  = -1
    10 = 4
    13 = -1
    15 = -5
    16 = -2
    21 = 1
    23 = 0
    26 = -2
    29 = 0
    37 = -3
    38 = -3
    40 = 1
    44 = -2
    51 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
